# Update Name of Algo
# Applies updated numeric results to columns D and E on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "E7"  = 16.609
    "D8"  = -8.559000000000001
    "D10" = -8.403
    "D12" = -7.290000000000001
    "E15" = 16.367
    "D18" = -8.318999999999999
    "E18" = 16.353
    "E20" = 16.396
    "E29" = 17.05
    "E30" = 16.578
    "E31" = 16.568
    "D37" = -8.081999999999999
    "E40" = 16.627
    "E50" = 16.249
    "D55" = -8.319000000000001
    "D68" = -7.188
    "E68" = 17.777
    "E76" = 16.558
    "D77" = -7.840000000000001
    "D78" = -8.08
    "D81" = -7.536
    "D82" = -8.154
    "E87" = 16.38
    "E88" = 16.208
    "E96" = 16.325
    "E98" = 16.299
    "E101" = 16.625
    "E102" = 16.617
}

foreach ($key in $updates.Keys) {
    $ws.Range($key).Value = $updates[$key]
}
